$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update Moisture Content of Feed Wood (%) values from 20 to 0 for rows 3-9
$ws.Range("G3:G9").Value = 0

# Update the selection on the Data sheet to G10
$ws.Activate()
$ws.Range("G10").Select()
